$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 3 new task rows before the old row 29 -------------------------
# (old rows 29-44 shift down to 32-47)
$ws.Range("A29:A31").EntireRow.Insert()
$ws.Rows.Item(29).RowHeight = 30
$ws.Rows.Item(30).RowHeight = 30
$ws.Rows.Item(31).RowHeight = 30

# New row 29: "3.13 refactor previous work"
$ws.Range("B29").Value = "3.13 refactor previous work"
$ws.Range("C29").Value = 20
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 20
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 0.8

# New row 30: "3.14 design correctness part"
$ws.Range("B30").Value = "3.14 design correctness part"
$ws.Range("C30").Value = 20
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 20
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 1

# New row 31: "3.15 correctness tutorial"
$ws.Range("B31").Value = "3.15 correctness tutorial"
$ws.Range("C31").Value = 21
$ws.Range("D31").Value = 1
$ws.Range("E31").Value = 21
$ws.Range("F31").Value = 1
$ws.Range("G31").ClearContents()

# Old row 29 ("3.7 Test algorithm correctness") is now row 32 - plan start
# moved from week 18 to week 21, and an actual duration of 2 was logged.
$ws.Range("C32").Value = 21
$ws.Range("F32").Value = 2

# --- Re-point the conditional formatting ranges to the new row count ------
# NB: Range(...).FormatConditions matches by intersection, not exact sqref,
# so B45:BO45 must be re-pointed to B48:BO48 *before* H5:AI44 grows down to
# H5:AI47 (which would otherwise start overlapping row 45-47 of B45:BO45 and
# get swept up together with it).
$cf2 = $ws.Range("B45:BO45").FormatConditions
for ($i = 1; $i -le $cf2.Count(); $i++) {
    $cf2.Item($i).ModifyAppliesToRange($ws.Range("B48:BO48"))
}

$cf1 = $ws.Range("H5:AI44").FormatConditions
for ($i = 1; $i -le $cf1.Count(); $i++) {
    $cf1.Item($i).ModifyAppliesToRange($ws.Range("H5:AI47"))
}

# --- Update the on-screen selection to match the new layout ---------------
$ws.Activate()
$ws.Range("G31").Select()
